$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 163.25
$ws.Range("I12").Value = 163.25
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 163.25
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 6.75
$ws.Range("N12").ClearContents()  # was -369
$ws.Range("H113").Value = 3688.25
$ws.Range("J113").Value = 3715.1428
$ws.Range("L113").Value = 3715.1428
$ws.Range("N113").Value = -10223.1428
$ws.Range("H132").Value = 11502745
$ws.Range("I132").Value = 13339933
$ws.Range("K132").Value = 40019799
$ws.Range("M132").Value = -40017269
$ws.Range("H135").Value = 542.1
$ws.Range("I135").Value = 255.41176
$ws.Range("K135").Value = 2298.70584
$ws.Range("M135").Value = 236.2941599999999
$ws.Range("H138").Value = 1335.14
$ws.Range("I138").Value = 870.7143
$ws.Range("J138").Value = 1585.2153
$ws.Range("K138").Value = 2612.1429
$ws.Range("L138").Value = 4755.6459
$ws.Range("M138").Value = 2527.8571
$ws.Range("N138").Value = -15035.6459

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 25703.25
$ws.Range("I2").Value = 600
$ws.Range("J2").Value = 34071
$ws.Range("K2").Value = 600
$ws.Range("L2").Value = 34071
$ws.Range("M2").Value = -487
$ws.Range("N2").Value = -34297
$ws.Range("H19").Value = 3750
$ws.Range("I19").Value = 3500
$ws.Range("K19").Value = 3500
$ws.Range("M19").Value = -3271
$ws.Range("H32").Value = 3198.3967
$ws.Range("I32").Value = 3391.3542
$ws.Range("J32").Value = 2580.9333
$ws.Range("K32").Value = 3391.3542
$ws.Range("L32").Value = 2580.9333
$ws.Range("M32").Value = -3104.3542
$ws.Range("N32").Value = -3154.9333
$ws.Range("H45").Value = 1123.963
$ws.Range("I45").Value = 1177.8572
$ws.Range("K45").Value = 1177.8572
$ws.Range("M45").Value = -800.8571999999999
$ws.Range("H74").Value = 1126.4062
$ws.Range("J74").Value = 2044.8889
$ws.Range("L74").Value = 2044.8889
$ws.Range("N74").Value = -3792.8889
$ws.Range("H77").Value = 1126.4062
$ws.Range("J77").Value = 2044.8889
$ws.Range("L77").Value = 10224.4445
$ws.Range("N77").Value = -18960.4445
$ws.Range("H110").Value = 1166.65
$ws.Range("I110").Value = 647.13336
$ws.Range("J110").Value = 2725.2
$ws.Range("K110").Value = 647.13336
$ws.Range("L110").Value = 2725.2
$ws.Range("M110").Value = 1397.86664
$ws.Range("N110").Value = -6815.2
$ws.Range("H116").Value = 25703.25
$ws.Range("I116").Value = 600
$ws.Range("J116").Value = 34071
$ws.Range("K116").Value = 600
$ws.Range("L116").Value = 34071
$ws.Range("M116").Value = 1694
$ws.Range("N116").Value = -38659

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 25703.25
$ws.Range("I3").Value = 600
$ws.Range("J3").Value = 34071
$ws.Range("K3").Value = 600
$ws.Range("L3").Value = 34071
$ws.Range("M3").Value = -486
$ws.Range("N3").Value = -34299
$ws.Range("H99").Value = 23810840
$ws.Range("I99").Value = 31251334
$ws.Range("J99").Value = 1260
$ws.Range("K99").Value = 31251334
$ws.Range("L99").Value = 1260
$ws.Range("M99").Value = -31249836
$ws.Range("N99").Value = -4256
$ws.Range("H105").Value = 77686250
$ws.Range("I105").Value = 91810530
$ws.Range("K105").Value = 91810530
$ws.Range("M105").Value = -91808783
$ws.Range("H107").Value = 2093.2727
$ws.Range("I107").Value = 1725
$ws.Range("K107").Value = 1725
$ws.Range("M107").Value = 195

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 66667676
$ws.Range("I16").Value = 71429510
$ws.Range("K16").Value = 71429510
$ws.Range("M16").Value = -71429223
$ws.Range("H58").Value = 983.7959
$ws.Range("I58").Value = 782.55884
$ws.Range("J58").Value = 1439.9333
$ws.Range("K58").Value = 782.55884
$ws.Range("L58").Value = 1439.9333
$ws.Range("M58").Value = -579.55884
$ws.Range("N58").Value = -1845.9333
$ws.Range("H99").Value = 2017.3462
$ws.Range("I99").Value = 1713.8422
$ws.Range("J99").Value = 2841.1428
$ws.Range("K99").Value = 1713.8422
$ws.Range("L99").Value = 2841.1428
$ws.Range("M99").Value = -215.8422
$ws.Range("N99").Value = -5837.1428
$ws.Range("H113").Value = 66667676
$ws.Range("I113").Value = 71429510
$ws.Range("K113").Value = 71429510
$ws.Range("M113").Value = -71427340
$ws.Range("H126").Value = 2017.3462
$ws.Range("I126").Value = 1713.8422
$ws.Range("J126").Value = 2841.1428
$ws.Range("K126").Value = 5141.5266
$ws.Range("L126").Value = 8523.428400000001
$ws.Range("M126").Value = -2671.5266
$ws.Range("N126").Value = -13463.4284
$ws.Range("H132").Value = 5032.697
$ws.Range("I132").Value = 5807.6816
$ws.Range("K132").Value = 17423.0448
$ws.Range("M132").Value = -14893.0448
$ws.Range("H134").Value = 2354.9
$ws.Range("I134").Value = 2792.5
$ws.Range("J134").Value = 1698.5
$ws.Range("K134").Value = 8377.5
$ws.Range("L134").Value = 5095.5
$ws.Range("M134").Value = -5842.5
$ws.Range("N134").Value = -10165.5
$ws.Range("H136").Value = 983.7959
$ws.Range("I136").Value = 782.55884
$ws.Range("J136").Value = 1439.9333
$ws.Range("K136").Value = 2347.67652
$ws.Range("L136").Value = 4319.7999
$ws.Range("M136").Value = 202.32348
$ws.Range("N136").Value = -9419.7999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 866.6667
$ws.Range("J107").Value = 866.6667
$ws.Range("L107").Value = 2600.0001
$ws.Range("N107").Value = -6440.0001
$ws.Range("H109").Value = 55530.527
$ws.Range("I109").Value = 167530
$ws.Range("K109").Value = 502590
$ws.Range("M109").Value = -501550
$ws.Range("H131").Value = 20001282
$ws.Range("I131").Value = 100000490
$ws.Range("J131").Value = 1481.2
$ws.Range("K131").Value = 300001470
$ws.Range("L131").Value = 4443.6
$ws.Range("M131").Value = -299996430
$ws.Range("N131").Value = -14523.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2299.2
$ws.Range("I102").Value = 2349
$ws.Range("K102").Value = 2349
$ws.Range("M102").Value = -727
$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 1100
$ws.Range("K113").Value = 1100
$ws.Range("M113").Value = 1070
$ws.Range("H126").Value = 2007.2778
$ws.Range("I126").Value = 1935.9166
$ws.Range("J126").Value = 2150
$ws.Range("K126").Value = 5807.7498
$ws.Range("L126").Value = 6450
$ws.Range("M126").Value = -3337.7498
$ws.Range("N126").Value = -11390
$ws.Range("H132").Value = 2542.8823
$ws.Range("I132").Value = 2094.6924
$ws.Range("K132").Value = 6284.0772
$ws.Range("M132").Value = -3754.0772

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 24340
$ws.Range("J116").Value = 24340
$ws.Range("L116").Value = 24340
$ws.Range("N116").Value = -33518
$ws.Range("H132").Value = 23021.469
$ws.Range("I132").Value = 1486.0476
$ws.Range("J132").Value = 40415.46
$ws.Range("K132").Value = 4458.142800000001
$ws.Range("L132").Value = 121246.38
$ws.Range("M132").Value = -1928.142800000001
$ws.Range("N132").Value = -126306.38

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 468.5
$ws.Range("I100").Value = 424.66666
$ws.Range("K100").Value = 849.33332
$ws.Range("M100").Value = -308.33332
$ws.Range("H113").Value = 471.5
$ws.Range("I113").Value = 301.66666
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 904.9999799999999
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = 1265.00002
$ws.Range("N113").Value = -10340
$ws.Range("H126").Value = 40000730
$ws.Range("I126").Value = 52632390
$ws.Range("J126").Value = 462.33334
$ws.Range("K126").Value = 157897170
$ws.Range("L126").Value = 1387.00002
$ws.Range("M126").Value = -157894700
$ws.Range("N126").Value = -6327.000019999999
